# Apply the "Testing commit -a" edit to the single paragraph in Skills.docx.
#
# Target state (from the diff):
#   - The paragraph mark (pilcrow) gets an explicit language of en-US,
#     i.e. a <w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr> is added.
#   - The first run's text "skill" becomes "Testing commit -a" and that run's
#     rPr gains <w:lang w:val="en-US"/> (all other run formatting is kept).
#   - The second run (a lone, differently-sized space character) is removed
#     entirely.
#
# Rather than trying to coax the COM Range/Font/LanguageID setters into
# creating a standalone paragraph-mark <w:pPr><w:rPr> (this host always folds
# paragraph-mark-only formatting into the adjacent run), we rebuild the whole
# paragraph's OOXML directly and drop it in with Range.InsertXML, which
# replaces the contents of the target range with exactly the XML supplied.

$d = $word.ActiveDocument
$para = $d.Paragraphs(1)

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6A237EA1" w14:textId="724AB083" w:rsidR="00182625" w:rsidRDefault="00477C3A" w:rsidP="00E62203"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Cascadia Mono" w:hAnsi="Cascadia Mono" w:cs="Cascadia Mono"/><w:color w:val="000000" w:themeColor="text1"/><w:kern w:val="0"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>Testing commit -a</w:t></w:r></w:p>'

$null = $para.Range.InsertXML($newParaXml)

Write-Output "Paragraph replaced."
